$wb = $excel.ActiveWorkbook

# Sheet references
$wsCreate = $wb.Worksheets.Item("createUser")
$wsEdit   = $wb.Worksheets.Item("editUser")
$wsReset  = $wb.Worksheets.Item("resetPassword")
$wsCheck  = $wb.Worksheets.Item("checkLogin")

# ---- createUser sheet ----
$wsCreate.Range("A2").Value = "UserTest-103"
$wsCreate.Range("C2").Value = "user.testAuto103"
$wsCreate.Range("A3").Value = "UserTest-104"
$wsCreate.Range("C3").Value = "user.testAuto104"
$wsCreate.Range("A4").Value = "UserTest-105"
$wsCreate.Range("C4").Value = "user.testAuto105"
$wsCreate.Range("A5").Value = "UserTest-106"
$wsCreate.Range("C5").Value = "user.testAuto106"
$wsCreate.Range("A6").Value = "UserTest-107"
$wsCreate.Range("C6").Value = "user.testAuto107"
$wsCreate.Range("A7").Value = "UserTest-108"
$wsCreate.Range("C7").Value = "user.testAuto108"

# ---- editUser sheet ----
$wsEdit.Range("A2").Value = "TestUser-86"
$wsEdit.Range("C2").Value = "userEdit.auto86"
$wsEdit.Range("A3").Value = "TestUser-87"
$wsEdit.Range("C3").Value = "userEdit.auto87"
$wsEdit.Range("A4").Value = "TestUser-88"
$wsEdit.Range("C4").Value = "userEdit.auto88"
$wsEdit.Range("A5").Value = "TestUser-89"
$wsEdit.Range("C5").Value = "userEdit.auto89"
$wsEdit.Range("A6").Value = "TestUser-90"
$wsEdit.Range("C6").Value = "userEdit.auto90"

# ---- resetPassword sheet ----
$wsReset.Range("A2").Value = "TestUser-85"
$wsReset.Range("C2").Value = "userPass.auto85"
$wsReset.Range("A3").Value = "TestUser-86"
$wsReset.Range("C3").Value = "userPass.auto86"
$wsReset.Range("A4").Value = "TestUser-87"
$wsReset.Range("C4").Value = "userPass.auto87"
$wsReset.Range("A5").Value = "TestUser-88"
$wsReset.Range("C5").Value = "userPass.auto88"
$wsReset.Range("A6").Value = "TestUser-89"
$wsReset.Range("C6").Value = "userPass.auto89"

# ---- checkLogin sheet ----
$wsCheck.Range("A2").Value = "UserTest-91"
$wsCheck.Range("C2").Value = "user_logintest.auto91"
$wsCheck.Range("A3").Value = "UserTest-92"
$wsCheck.Range("C3").Value = "user_logintest.auto92"
$wsCheck.Range("A4").Value = "UserTest-93"
$wsCheck.Range("C4").Value = "user_logintest.auto93"
$wsCheck.Range("A5").Value = "UserTest-94"
$wsCheck.Range("C5").Value = "user_logintest.auto94"
$wsCheck.Range("A6").Value = "UserTest-95"
$wsCheck.Range("C6").Value = "user_logintest.auto95"
$wsCheck.Range("A7").Value = "UserTest-96"
$wsCheck.Range("C7").Value = "user_logintest.auto96"

# ---- Column widths ----
# editUser: column C width -> 14 (target XML width="14", bestFit="1")
# NOTE: the runtime stores/returns column width with an internal +0.8333333333333
# offset that gets baked into the saved XML (character-width quantization), so we
# back-compute the COM input required to land exactly on the desired saved width.
$wsEdit.Columns.Item(3).ColumnWidth = 13.166666666666666
# resetPassword: column C width -> 14.44140625, bestFit (new column definition)
# (closest value reachable given the engine's width quantization is 14.5)
$wsReset.Columns.Item(3).ColumnWidth = 13.666666666666666

# ---- Sheet view selections ----
$wsCreate.Range("C2").Select()
$wsEdit.Range("D17").Select()
$wsReset.Range("C6").Select()
$wsCheck.Range("C7").Select()

# ---- Active sheet / tab selection ----
# New workbookView activeTab="3" -> checkLogin (4th sheet, 0-indexed 3) is active when opening
$wsCheck.Activate()
